$d = $word.ActiveDocument

# Locate the paragraph ending in "Ou seja, o JS faz a interação do cliente com o site."
# That paragraph is immediately followed (in the original document) by a single trailing
# empty paragraph. We append the new block of paragraphs right after it.
$anchorIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -eq "Ou seja, o JS faz a interação do cliente com o site.`r") {
        $anchorIndex = $i
    }
}
if ($anchorIndex -eq 0) {
    throw "Anchor paragraph 'Ou seja, o JS faz a interação do cliente com o site.' not found"
}

$idx = $anchorIndex

# 1) blank paragraph
$d.Paragraphs($idx).Range.InsertParagraphAfter()
$idx = $idx + 1

# 2) "Vamos a um exemplo de programação WEB" + ":"
$d.Paragraphs($idx).Range.InsertParagraphAfter()
$idx = $idx + 1
$d.Paragraphs($idx).Range.InsertAfter("Vamos a um exemplo de programação WEB")
$d.Paragraphs($idx).Range.InsertAfter(":")

# 3) HTML line
$d.Paragraphs($idx).Range.InsertParagraphAfter()
$idx = $idx + 1
$d.Paragraphs($idx).Range.InsertAfter("HTML – É toda escrita do site e suas semânticas")

# 4) CSS line
$d.Paragraphs($idx).Range.InsertParagraphAfter()
$idx = $idx + 1
$d.Paragraphs($idx).Range.InsertAfter("CSS – É todo estilo do site e sua beleza")

# 5) Java Script line
$d.Paragraphs($idx).Range.InsertParagraphAfter()
$idx = $idx + 1
$d.Paragraphs($idx).Range.InsertAfter("Java Script – Interatividade com o site, como quando você passa o mouse por cima de alguma imagem e aparece a descrição dela")

# 6) blank paragraph
$d.Paragraphs($idx).Range.InsertParagraphAfter()
$idx = $idx + 1

# 7) Observação paragraph
$d.Paragraphs($idx).Range.InsertParagraphAfter()
$idx = $idx + 1
$d.Paragraphs($idx).Range.InsertAfter("Observação: HTML e CSS são linguagens de construção de sites e não linguagens de programação, já o Java Script é uma linguagem de programação.")

Write-Output "Inserted block after paragraph $anchorIndex; last new paragraph index $idx; total paragraphs now $($d.Paragraphs.Count)"
